$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3830520358464513
$ws.Range("C2").Value = 0.8947454165326297
$ws.Range("D2").Value = 0.6327228580167067
$ws.Range("E2").Value = 0.4240862509183345

$ws.Range("B3").Value = 0.7596543070440361
$ws.Range("C3").Value = 0.6465080687461903
$ws.Range("D3").Value = 0.59483642934008
$ws.Range("E3").Value = 0.2693498587946093

$ws.Range("B4").Value = 0.4369891931017392
$ws.Range("C4").Value = 0.2750477894891764
$ws.Range("D4").Value = 0.5590702820641964
$ws.Range("E4").Value = 0.2384682983040076

$ws.Range("B5").Value = 0.761972593927858
$ws.Range("C5").Value = 0.9404414338443914
$ws.Range("D5").Value = 0.1252523231592332
$ws.Range("E5").Value = 0.6136490322132813

$ws.Range("B6").Value = 0.2120530015793479
$ws.Range("C6").Value = 0.945671954701841
$ws.Range("D6").Value = 0.6392771840427417
$ws.Range("E6").Value = 0.8740041794208206

$ws.Range("B7").Value = 0.9141177053444617
$ws.Range("C7").Value = 0.8362416142071636
$ws.Range("D7").Value = 0.2074197608251617
$ws.Range("E7").Value = 0.2823075770737018
